$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the old "DiD" column (I), shifting old I:L headers/data -> M:P
$ws.Range("I1:L1").EntireColumn.Insert()

# Header row: A1:H1 unchanged, I1:L1 are new headers, M1:P1 are the former I1:L1 headers
$ws.Range("A1").Value = 'KPI'
$ws.Range("B1").Value = 'Pre_Test_Mean'
$ws.Range("C1").Value = 'Pre_Control_Mean'
$ws.Range("D1").Value = 'Post_Test_Mean'
$ws.Range("E1").Value = 'Post_Control_Mean'
$ws.Range("F1").Value = 'Change_Test_Mean'
$ws.Range("G1").Value = 'Change_Control_Mean'
$ws.Range("H1").Value = 'Diff_in_Change'
$ws.Range("I1").Value = '%Change_Test'
$ws.Range("J1").Value = '%Change_Control'
$ws.Range("K1").Value = '%Change_Diff'
$ws.Range("L1").Value = 'Direction'
$ws.Range("M1").Value = 'DiD'
$ws.Range("N1").Value = 'tstat'
$ws.Range("O1").Value = 'pval'
$ws.Range("P1").Value = 'Significant'

# Data rows 2-8 (recomputed KPI summary stats)
# Row 2: LPE
$ws.Range("A2").Value = 'LPE'
$ws.Range("B2").Value = 0.2716196355085244
$ws.Range("C2").Value = 0.2486700312902269
$ws.Range("D2").Value = 0.2508465608465608
$ws.Range("E2").Value = 0.2201635211197646
$ws.Range("F2").Value = -0.02077307466196356
$ws.Range("G2").Value = -0.02850651017046224
$ws.Range("H2").Value = 0.007733435508498688
$ws.Range("I2").Value = -7.647854553325054
$ws.Range("J2").Value = -11.4635889264806
$ws.Range("K2").Value = 3.815734373155545
$ws.Range("L2").Value = '↑ Better'
$ws.Range("M2").Value = 0.007733435508498698
$ws.Range("N2").Value = 0.1930562810330638
$ws.Range("O2").Value = 0.8498495911476762
$ws.Range("P2").Value = $False

# Row 3: avg_loan_size
$ws.Range("A3").Value = 'avg_loan_size'
$ws.Range("B3").Value = 7449.012972222223
$ws.Range("C3").Value = 7287.316710643691
$ws.Range("D3").Value = 6773.931223544973
$ws.Range("E3").Value = 6825.928893044526
$ws.Range("F3").Value = -1318.439095734127
$ws.Range("G3").Value = -461.3878175991657
$ws.Range("H3").Value = -857.0512781349619
$ws.Range("I3").Value = -9.062700671816067
$ws.Range("J3").Value = -6.331381438730009
$ws.Range("K3").Value = -2.731319233086058
$ws.Range("L3").Value = '↓ Worse'
$ws.Range("M3").Value = -213.6939310780845
$ws.Range("N3").Value = -1.050015317374909
$ws.Range("O3").Value = 0.317159510780542
$ws.Range("P3").Value = $False

# Row 4: dq29_pot30_payment_rate_$_up_to_day
$ws.Range("A4").Value = 'dq29_pot30_payment_rate_$_up_to_day'
$ws.Range("B4").Value = 0.2693797627988484
$ws.Range("C4").Value = 0.251513158863541
$ws.Range("D4").Value = 0.2874381755342127
$ws.Range("E4").Value = 0.2625666299828348
$ws.Range("F4").Value = 0.01805841273536432
$ws.Range("G4").Value = 0.01105347111929382
$ws.Range("H4").Value = 0.007004941616070504
$ws.Range("I4").Value = 6.703700585277036
$ws.Range("J4").Value = 4.394788395660367
$ws.Range("K4").Value = 2.308912189616669
$ws.Range("L4").Value = '↑ Better'
$ws.Range("M4").Value = 0.00700494161607057
$ws.Range("N4").Value = 1.17674528642536
$ws.Range("O4").Value = 0.2667296177651306
$ws.Range("P4").Value = $False

# Row 5: dq29_pot30_payment_rate_unit_per_day
$ws.Range("A5").Value = 'dq29_pot30_payment_rate_unit_per_day'
$ws.Range("B5").Value = 0.007372493582211477
$ws.Range("C5").Value = 0.008094526316713242
$ws.Range("D5").Value = 0.00682756008862527
$ws.Range("E5").Value = 0.007112913700804932
$ws.Range("F5").Value = -0.0005449334935862068
$ws.Range("G5").Value = -0.000981612615908309
$ws.Range("H5").Value = 0.0004366791223221021
$ws.Range("I5").Value = -7.391440731817451
$ws.Range("J5").Value = -12.12686916443174
$ws.Range("K5").Value = 4.735428432614285
$ws.Range("L5").Value = '↑ Better'
$ws.Range("M5").Value = 0.0004366791223221035
$ws.Range("N5").Value = 0.3290343122917556
$ws.Range("O5").Value = 0.7488612092866946
$ws.Range("P5").Value = $False

# Row 6: dq29_pot30_payment_rate_unit_up_to_day
$ws.Range("A6").Value = 'dq29_pot30_payment_rate_unit_up_to_day'
$ws.Range("B6").Value = 0.7290152249144803
$ws.Range("C6").Value = 0.682195370172997
$ws.Range("D6").Value = 0.3866628639926504
$ws.Range("E6").Value = 0.4290904643047829
$ws.Range("F6").Value = -0.34235236092183
$ws.Range("G6").Value = -0.253104905868214
$ws.Range("H6").Value = -0.08924745505361598
$ws.Range("I6").Value = -46.9609343154652
$ws.Range("J6").Value = -37.1015279397205
$ws.Range("K6").Value = -9.859406375744697
$ws.Range("L6").Value = '↓ Worse'
$ws.Range("M6").Value = -0.08924745505361581
$ws.Range("N6").Value = -3.856119468263016
$ws.Range("O6").Value = 0.003976677556247839
$ws.Range("P6").Value = $True

# Row 7: dq30_pct_$
$ws.Range("A7").Value = 'dq30_pct_$'
$ws.Range("B7").Value = 0.6924078466606944
$ws.Range("C7").Value = 0.7346892312631622
$ws.Range("D7").Value = 0.6753476016538416
$ws.Range("E7").Value = 0.7427374434574096
$ws.Range("F7").Value = -0.01706024500685287
$ws.Range("G7").Value = 0.008048212194247234
$ws.Range("H7").Value = -0.0251084572011001
$ws.Range("I7").Value = -2.463901165928444
$ws.Range("J7").Value = 1.095458032018508
$ws.Range("K7").Value = -3.559359197946951
$ws.Range("L7").Value = '↓ Worse'
$ws.Range("M7").Value = -0.02510845720110022
$ws.Range("N7").Value = -1.919468142800438
$ws.Range("O7").Value = 0.08633284165486863
$ws.Range("P7").Value = $False

# Row 8: dq30_pct_unit
$ws.Range("A8").Value = 'dq30_pct_unit'
$ws.Range("B8").Value = 0.05369644062367666
$ws.Range("C8").Value = 0.05033816976596275
$ws.Range("D8").Value = 0.03388713078061231
$ws.Range("E8").Value = 0.03331890603899778
$ws.Range("F8").Value = -0.01980930984306436
$ws.Range("G8").Value = -0.01701926372696497
$ws.Range("H8").Value = -0.002790046116099388
$ws.Range("I8").Value = -36.89129039649927
$ws.Range("J8").Value = -33.80985801846318
$ws.Range("K8").Value = -3.081432378036084
$ws.Range("L8").Value = '↓ Worse'
$ws.Range("M8").Value = -0.002790046116099391
$ws.Range("N8").Value = -1.117266812063672
$ws.Range("O8").Value = 0.2883794507195882
$ws.Range("P8").Value = $False

